# Added default currency formatting.
#
# This edit:
#  1. Adds a new "currency" column (H) to the boolean/numeric/string type
#     transform table (rows 1-6), describing whether each type can be
#     transformed to/from :currency.
#  2. Renames the ":utc_datetime" column header (E8) in the date/time type
#     transform table to ":datetime".
#  3. Inserts a new row for the ":datetime" (-> "DateTime") transform
#     target, between the existing ":naive_datetime" (NaiveDateTime) and
#     Integer rows, pushing the rest of that table down by one row.
#  4. Adds a clarifying note for the (now shifted) Integer row explaining
#     that the integer is assumed to be a unix timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Date/time transform table -------------------------------------------

# Header: ":utc_datetime" -> ":datetime"
$ws.Range("E8").Value2 = ":datetime"

# Make room for the new "DateTime" row by inserting a row above the old
# "Integer" row (row 11) - this shifts the old rows 11-13 (Integer, Time,
# String) down to 12-14, and the trailing note row from 16 to 17.
$ws.Rows("11").Insert()

# Populate the newly inserted row 11 with the DateTime transform data.
$ws.Range("A11").Value2 = "DateTime"
$ws.Range("C11").Value2 = "T -> 00:00:00"
$ws.Range("D11").Value2 = "Y"
$ws.Range("F11").Value2 = "00:00:00"
$ws.Range("G11").Value2 = '"YYYY-MM-DD hh:mm:ssZ"'

# Clarifying note alongside the (shifted) Integer row.
$ws.Range("I12").Value2 = "Assumes integer is unix_time"

# --- Boolean/numeric/string transform table -------------------------------

# New ":currency" column, formatted like the existing columns.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value2 = ":currency"
$ws.Range("H2").Value2 = "N"
$ws.Range("H3").Value2 = "Y"
$ws.Range("H4").Value2 = "Y"
$ws.Range("H5").Value2 = "Y"
$ws.Range("H6").Value2 = "Y (S->D->C)"

# Restore the selection to where the edit left off.
$null = $ws.Range("G18").Select()
